$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 137, shifting rows 137:149 down to 139:151
$ws.Range("A137:R138").EntireRow.Insert()

# Row 137: Coliflor "Primera" entry for Vega Monumental Concepción, new date 44474
$ws.Range("A137").Value = 11
$ws.Range("B137").Value = "Vega Monumental Concepción"
$ws.Range("C137").Value = "Bíobío"
$ws.Range("D137").Value = 44474
$ws.Range("D137").NumberFormat = $ws.Range("D139").NumberFormat
$ws.Range("E137").Value = 8
$ws.Range("F137").Value = 100112008
$ws.Range("G137").Value = "Coliflor"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 600
$ws.Range("L137").Value = 700
$ws.Range("M137").Value = 650
$ws.Range("N137").Value = "$/unidad"
$ws.Range("O137").Value = "Región Metropolitana"
$ws.Range("P137").Value = 650
$ws.Range("Q137").Value = 1
$ws.Range("R137").Value = "Hortaliza"

# Row 138: Coliflor "Segunda" entry for Vega Monumental Concepción, new date 44474
$ws.Range("A138").Value = 11
$ws.Range("B138").Value = "Vega Monumental Concepción"
$ws.Range("C138").Value = "Bíobío"
$ws.Range("D138").Value = 44474
$ws.Range("D138").NumberFormat = $ws.Range("D139").NumberFormat
$ws.Range("E138").Value = 8
$ws.Range("F138").Value = 100112008
$ws.Range("G138").Value = "Coliflor"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Segunda"
$ws.Range("J138").Value = 500
$ws.Range("K138").Value = 500
$ws.Range("L138").Value = 500
$ws.Range("M138").Value = 500
$ws.Range("N138").Value = "$/unidad"
$ws.Range("O138").Value = "Región Metropolitana"
$ws.Range("P138").Value = 500
$ws.Range("Q138").Value = 1
$ws.Range("R138").Value = "Hortaliza"
